# implement init filtering logic
# Translate the Arabic ledger sheet to English, rename the sheet/tab,
# and restore the "A2" selection that Excel persisted on last save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet (tab) name: "ورقة1" -> "sheet1" ---
$ws.Name = "sheet1"

# --- header block (rows 2-5) ---
$ws.Range("A2").Value = "Cashier Number"
$ws.Range("D2").Value = "Exhibit"
$ws.Range("E2").Value = "Currency"
$ws.Range("F3").Value = "Saudi Riyal"

$ws.Range("A4").Value = "Date"
$ws.Range("B4").Value = "Type"
$ws.Range("C4").Value = "Number"
$ws.Range("D4").Value = "Description"
$ws.Range("E4").Value = "Ref Number"
$ws.Range("F4").Value = "Credit"
$ws.Range("G4").Value = "Debit"

$ws.Range("D5").Value = "Open Balance"

# --- transaction rows ---
$ws.Range("B7").Value = "Journal Entry"
$ws.Range("D7").Value = "Cash purchase invoices"

$ws.Range("B8").Value = "Payment Voucher"
$ws.Range("B9").Value = "Payment Voucher"
$ws.Range("B10").Value = "Payment Voucher"

$ws.Range("B11").Value = "Invoice"
$ws.Range("D11").Value = "Sales"
$ws.Range("B12").Value = "Invoice"
$ws.Range("D12").Value = "Sales"
$ws.Range("B13").Value = "Invoice"
$ws.Range("D13").Value = "Sales"
$ws.Range("B14").Value = "Invoice"
$ws.Range("D14").Value = "Sales"
$ws.Range("B15").Value = "Invoice"
$ws.Range("D15").Value = "Sales"
$ws.Range("B16").Value = "Invoice"
$ws.Range("D16").Value = "Sales"
$ws.Range("B17").Value = "Invoice"
$ws.Range("D17").Value = "Sales"
$ws.Range("B18").Value = "Invoice"
$ws.Range("D18").Value = "Sales"
$ws.Range("B19").Value = "Invoice"
$ws.Range("D19").Value = "Sales"
$ws.Range("B20").Value = "Invoice"
$ws.Range("D20").Value = "Sales"
$ws.Range("B21").Value = "Invoice"
$ws.Range("D21").Value = "Sales"
$ws.Range("B22").Value = "Invoice"
$ws.Range("D22").Value = "Sales"
$ws.Range("B23").Value = "Invoice"
$ws.Range("D23").Value = "Sales"
$ws.Range("B24").Value = "Invoice"
$ws.Range("D24").Value = "Sales"
$ws.Range("B25").Value = "Invoice"
$ws.Range("D25").Value = "Sales"
$ws.Range("B26").Value = "Invoice"
$ws.Range("D26").Value = "Sales"
$ws.Range("B27").Value = "Invoice"
$ws.Range("D27").Value = "Sales"

$ws.Range("B28").Value = "Return"
$ws.Range("D28").Value = "Return for Invoice No. 3625"
$ws.Range("B29").Value = "Return"
$ws.Range("D29").Value = "Return"
$ws.Range("B30").Value = "Return"
$ws.Range("D30").Value = "Return for Invoice No. 3599"
$ws.Range("B31").Value = "Return"
$ws.Range("D31").Value = "Return for Invoice No. 3631"

$ws.Range("B32").Value = "Purchase"
$ws.Range("B33").Value = "Purchase"

$ws.Range("D34").Value = "Total"

# --- font: Arial -> Calibri across the whole sheet ---
$ws.Range("A1:G34").Font.Name = "Calibri"

# --- restore the persisted selection (A2) ---
$ws.Range("A2").Select()
